$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H98").Value = 1140.5416
$ws_ALC.Range("I98").Value = 1051.6666
$ws_ALC.Range("K98").Value = 1051.6666
$ws_ALC.Range("M98").Value = 446.3334

$ws_ALC.Range("H122").Value = 1140.5416
$ws_ALC.Range("I122").Value = 1051.6666
$ws_ALC.Range("K122").Value = 3154.9998
$ws_ALC.Range("M122").Value = -704.9998000000001

$ws_ALC.Range("H132").Value = 2926.1462
$ws_ALC.Range("I132").Value = 2877.2646
$ws_ALC.Range("J132").Value = 3163.5715
$ws_ALC.Range("K132").Value = 8631.793799999999
$ws_ALC.Range("L132").Value = 9490.7145
$ws_ALC.Range("M132").Value = -6101.793799999999
$ws_ALC.Range("N132").Value = -14550.7145

$ws_ALC.Range("H137").Value = 1063.4
$ws_ALC.Range("I137").Value = 302
$ws_ALC.Range("J137").Value = 1253.75
$ws_ALC.Range("K137").Value = 906
$ws_ALC.Range("L137").Value = 3761.25
$ws_ALC.Range("M137").Value = 1644
$ws_ALC.Range("N137").Value = -8861.25

$ws_ALC.Range("H138").Value = 4499.143
$ws_ALC.Range("I138").Value = 3098.8
$ws_ALC.Range("J138").Value = 8000
$ws_ALC.Range("K138").Value = 9296.400000000001
$ws_ALC.Range("L138").Value = 24000
$ws_ALC.Range("M138").Value = -4156.400000000001
$ws_ALC.Range("N138").Value = -34280

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 3374.0686
$ws_ARM.Range("I32").Value = 3241.446
$ws_ARM.Range("K32").Value = 3241.446
$ws_ARM.Range("M32").Value = -2954.446

$ws_ARM.Range("H74").Value = 1960.7
$ws_ARM.Range("I74").Value = 1762.3572
$ws_ARM.Range("J74").Value = 2423.5
$ws_ARM.Range("K74").Value = 1762.3572
$ws_ARM.Range("L74").Value = 2423.5
$ws_ARM.Range("M74").Value = -888.3571999999999
$ws_ARM.Range("N74").Value = -4171.5

$ws_ARM.Range("H77").Value = 1960.7
$ws_ARM.Range("I77").Value = 1762.3572
$ws_ARM.Range("J77").Value = 2423.5
$ws_ARM.Range("K77").Value = 8811.786
$ws_ARM.Range("L77").Value = 12117.5
$ws_ARM.Range("M77").Value = -4443.786
$ws_ARM.Range("N77").Value = -20853.5

$ws_ARM.Range("H122").Value = 1686.52
$ws_ARM.Range("I122").Value = 1782.591
$ws_ARM.Range("K122").Value = 5347.772999999999
$ws_ARM.Range("M122").Value = -2897.772999999999

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H94").Value = 7937205.5
$ws_BSM.Range("I94").Value = 10417165
$ws_BSM.Range("J94").Value = 1336.4
$ws_BSM.Range("K94").Value = 10417165
$ws_BSM.Range("L94").Value = 1336.4
$ws_BSM.Range("M94").Value = -10416714
$ws_BSM.Range("N94").Value = -2238.4

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H7").Value = 192.64706
$ws_CRP.Range("I7").Value = 157.5
$ws_CRP.Range("J7").Value = 242.85715
$ws_CRP.Range("K7").Value = 157.5
$ws_CRP.Range("L7").Value = 242.85715
$ws_CRP.Range("M7").Value = -44.5
$ws_CRP.Range("N7").Value = -468.85715

$ws_CRP.Range("H31").Value = 7806.4424
$ws_CRP.Range("I31").Value = 2589.7222
$ws_CRP.Range("J31").Value = 10568.235
$ws_CRP.Range("K31").Value = 2589.7222
$ws_CRP.Range("L31").Value = 10568.235
$ws_CRP.Range("M31").Value = -2294.7222
$ws_CRP.Range("N31").Value = -11158.235

$ws_CRP.Range("H34").Value = 7806.4424
$ws_CRP.Range("I34").Value = 2589.7222
$ws_CRP.Range("J34").Value = 10568.235
$ws_CRP.Range("K34").Value = 2589.7222
$ws_CRP.Range("L34").Value = 10568.235
$ws_CRP.Range("M34").Value = -2387.7222
$ws_CRP.Range("N34").Value = -10972.235

$ws_CRP.Range("H50").Value = 46665.668
$ws_CRP.Range("I50").Value = 40000
$ws_CRP.Range("J50").Value = 49998.5
$ws_CRP.Range("K50").Value = 40000
$ws_CRP.Range("L50").Value = 49998.5
$ws_CRP.Range("M50").Value = -39375
$ws_CRP.Range("N50").Value = -51248.5

$ws_CRP.Range("H60").Value = 38986
$ws_CRP.Range("J60").Value = 38986
$ws_CRP.Range("L60").Value = 38986
$ws_CRP.Range("N60").Value = -40008

$ws_CRP.Range("H94").Value = 1600.8334
$ws_CRP.Range("I94").Value = 1719.5454
$ws_CRP.Range("K94").Value = 1719.5454
$ws_CRP.Range("M94").Value = -1268.5454

$ws_CRP.Range("H122").Value = 4110.8887
$ws_CRP.Range("I122").Value = 3796.524
$ws_CRP.Range("K122").Value = 11389.572
$ws_CRP.Range("M122").Value = -8939.572

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H56").Value = 13663.117
$ws_CUL.Range("I56").Value = 13663.117
$ws_CUL.Range("K56").Value = 13663.117
$ws_CUL.Range("M56").Value = -13133.117

$ws_CUL.Range("H68").Value = 3309.4285
$ws_CUL.Range("I68").Value = 2834
$ws_CUL.Range("J68").Value = 3439.0908
$ws_CUL.Range("K68").Value = 8502
$ws_CUL.Range("L68").Value = 10317.2724
$ws_CUL.Range("M68").Value = -7691
$ws_CUL.Range("N68").Value = -11939.2724

$ws_CUL.Range("H71").Value = 3309.4285
$ws_CUL.Range("I71").Value = 2834
$ws_CUL.Range("J71").Value = 3439.0908
$ws_CUL.Range("K71").Value = 25506
$ws_CUL.Range("L71").Value = 30951.8172
$ws_CUL.Range("M71").Value = -21450
$ws_CUL.Range("N71").Value = -39063.8172

$ws_CUL.Range("H80").Value = 25628.143
$ws_CUL.Range("I80").Value = 35499.668
$ws_CUL.Range("J80").Value = 18224.5
$ws_CUL.Range("K80").Value = 106499.004
$ws_CUL.Range("L80").Value = 54673.5
$ws_CUL.Range("M80").Value = -105563.004
$ws_CUL.Range("N80").Value = -56545.5

$ws_CUL.Range("H83").Value = 25628.143
$ws_CUL.Range("I83").Value = 35499.668
$ws_CUL.Range("J83").Value = 18224.5
$ws_CUL.Range("K83").Value = 319497.012
$ws_CUL.Range("L83").Value = 164020.5
$ws_CUL.Range("M83").Value = -314817.012
$ws_CUL.Range("N83").Value = -173380.5

$ws_CUL.Range("H130").Value = 0
$ws_CUL.Range("I130").Value = 0
$ws_CUL.Range("K130").Value = 0
$ws_CUL.Range("M130").Value = ""

$ws_CUL.Range("H131").Value = 2902.25
$ws_CUL.Range("I131").Value = 3343.6
$ws_CUL.Range("J131").Value = 2166.6667
$ws_CUL.Range("K131").Value = 10030.8
$ws_CUL.Range("L131").Value = 6500.000100000001
$ws_CUL.Range("M131").Value = -4990.799999999999
$ws_CUL.Range("N131").Value = -16580.0001

$ws_CUL.Range("H133").Value = 4900
$ws_CUL.Range("I133").Value = 4900
$ws_CUL.Range("K133").Value = 14700
$ws_CUL.Range("M133").Value = -9640

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H102").Value = 10420166
$ws_GSM.Range("I102").Value = 12198407
$ws_GSM.Range("K102").Value = 12198407
$ws_GSM.Range("M102").Value = -12196785

$ws_GSM.Range("H126").Value = 2988.1155
$ws_GSM.Range("I126").Value = 2907.84
$ws_GSM.Range("J126").Value = 4995
$ws_GSM.Range("K126").Value = 8723.52
$ws_GSM.Range("L126").Value = 14985
$ws_GSM.Range("M126").Value = -6253.52
$ws_GSM.Range("N126").Value = -19925

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H100").Value = 2552.3667
$ws_LTW.Range("I100").Value = 2356.4285
$ws_LTW.Range("J100").Value = 2723.8125
$ws_LTW.Range("K100").Value = 2356.4285
$ws_LTW.Range("L100").Value = 2723.8125
$ws_LTW.Range("M100").Value = -1815.4285
$ws_LTW.Range("N100").Value = -3805.8125

$ws_LTW.Range("H136").Value = 6538676.5
$ws_LTW.Range("I136").Value = 2250.139
$ws_LTW.Range("J136").Value = 22226100
$ws_LTW.Range("K136").Value = 6750.417
$ws_LTW.Range("L136").Value = 66678300
$ws_LTW.Range("M136").Value = -4200.417
$ws_LTW.Range("N136").Value = -66683400

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H139").Value = 77344.14999999999
$ws_WVR.Range("J139").Value = 77344.14999999999
$ws_WVR.Range("L139").Value = 77344.14999999999
$ws_WVR.Range("N139").Value = -87624.14999999999

$ws_WVR.Range("H141").Value = 111456
$ws_WVR.Range("J141").Value = 111456
$ws_WVR.Range("L141").Value = 111456
$ws_WVR.Range("N141").Value = -121816
